$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the DocType ID in B3: "Request##...response..." should be "Response##...response..."
$ws.Range("B3").Value = "urn:eu:toop:ns:dataexchange-1p10::Response##urn:eu.toop.response.registeredorganization::1.10"

# The longer text now wraps to two lines in the wrap-text cell, so the row grows taller
$ws.Rows.Item(3).RowHeight = 29

# Move selection to B3 (matches resulting selection in the file)
$ws.Range("B3").Select()
